$wb = $excel.ActiveWorkbook

# --- Reorder the sheets: review_info first, hotel_info second ---
$hotelInfo  = $wb.Worksheets.Item("hotel_info")
$reviewInfo = $wb.Worksheets.Item("review_info")
$reviewInfo.Move($hotelInfo)

# --- Add a new "State" column to hotel_info (between Hotel_Name and City) ---
$ws = $wb.Worksheets.Item("hotel_info")
$ws.Columns.Item(3).Insert()
$ws.Cells.Item(1, 3).Value = "State"
$ws.Cells.Item(2, 3).Value = "Louisiana"
